$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("D2").Value = '28.418.77'
$ws.Range("E2").Value = '  +3.35%  '
$ws.Range("D3").Value = '1.867.36'
$ws.Range("E3").Value = '  +1.94%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.50'
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4709'
$ws.Range("E7").Value = '  +2.35%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3975'
$ws.Range("E8").Value = '  +3.62%  '
$ws.Range("E9").Value = '  +2.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08019'
$ws.Range("E10").Value = '  +1.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9993'
$ws.Range("E11").Value = '  +2.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.00'
$ws.Range("E12").Value = '  +4.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.036'
$ws.Range("E13").Value = '  +2.62%  '
$ws.Range("D14").Value = '1.858.30'
$ws.Range("E14").Value = '  +0.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.254'
$ws.Range("E15").Value = '  +2.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.53'
$ws.Range("E16").Value = '  +2.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -0.36%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06647'
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001037'
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.53'
$ws.Range("E20").Value = '  +1.26%  '
$ws.Range("E21").Value = '  -0.33%  '
$ws.Range("D22").Value = '28.439.85'
$ws.Range("E22").Value = '  +3.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.476'
$ws.Range("E23").Value = '  +2.49%  '
$ws.Range("E24").Value = '  +2.26%  '
$ws.Range("E25").Value = '  -1.48%  '
$ws.Range("D26").Value = '2.083.92'
$ws.Range("E26").Value = '  +0.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.55'
$ws.Range("E27").Value = '  +2.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.75'
$ws.Range("E28").Value = '  +1.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.118'
$ws.Range("E29").Value = '  +2.65%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.480'
$ws.Range("E30").Value = '  +4.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '119.61'
$ws.Range("E31").Value = '  +1.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9637'
$ws.Range("E32").Value = '  +1.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09519'
$ws.Range("E33").Value = '  +2.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.593'
$ws.Range("E34").Value = '  +0.81%  '
$ws.Range("E35").Value = '  +4.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.361'
$ws.Range("E36").Value = '  +2.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06111'
$ws.Range("E37").Value = '  +2.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02251'
$ws.Range("E38").Value = '  +1.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.317'
$ws.Range("E39").Value = '  +3.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.180'
$ws.Range("E40").Value = '  +2.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5941'
$ws.Range("E41").Value = '  +2.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  -0.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1878'
$ws.Range("E43").Value = '  +1.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.32'
$ws.Range("E44").Value = '  +2.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.276'
$ws.Range("E45").Value = '  +0.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5569'
$ws.Range("E46").Value = '  +1.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.12'
$ws.Range("E47").Value = '  +1.40%  '
$ws.Range("E48").Value = '  +4.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07235'
$ws.Range("E49").Value = '  +8.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.068'
$ws.Range("E50").Value = '  +13.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '111.97'
